# Fix issues with cancel ticket:
#  1. Details sheet: add row 8 ("ram" / "123") below the existing data.
#  2. Add a new "walle" sheet (movie listing) with the same shape as the
#     existing Pushpa / Krish sheets.

# Helper: write a value that should be stored as TEXT even when it looks
# like a number (Excel would otherwise auto-detect it as numeric). A
# leading apostrophe is exactly what a user typing into the grid would do
# to force text, and is what the COM layer honours for quote-prefixed
# input.
function Set-TextValue($range, $value) {
    if ($value -match '^-?[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Details sheet (sheet1) -- append row 8
# ---------------------------------------------------------------------
$details = $wb.Worksheets.Item("Details")
Set-TextValue $details.Cells.Item(8, 1) "ram"
Set-TextValue $details.Cells.Item(8, 2) "123"

# ---------------------------------------------------------------------
# 2. New sheet "walle" -- inserted after the last existing sheet (Krish)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$walle = $wb.Worksheets.Add($null, $lastSheet)
$walle.Name = "walle"

$headers = @("title", "genre", "length", "cast", "director", "rating", "language", "timings", "shows_per_day", "firstshow", "interval_time", "gap", "capacity", "user_count")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $walle.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$rows = @(
    @("walle", "comedy", "120", "w&e", "ani", "5", "eng", "10:00-12:00", "3", "10", "15", "15", "200"),
    @("walle", "comedy", "120", "w&e", "ani", "5", "eng", "12:30-02:30", "3", "10", "15", "15", "200"),
    @("walle", "comedy", "120", "w&e", "ani", "5", "eng", "03:00-05:00", "3", "10", "15", "15", "200")
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowData = $rows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        Set-TextValue $walle.Cells.Item($excelRow, $c + 1) $rowData[$c]
    }
}

# user_count only has a value on the first data row, and it is a real
# number (not text).
$walle.Cells.Item(2, 14).Value = 1

# Adding the sheet shifted the active tab onto "walle"; restore the
# original active sheet (Pushpa) so the rest of the workbook's view state
# is left untouched.
$wb.Worksheets.Item("Pushpa").Activate()
